$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting existing rows 154:168 down to 155:169
$ws.Rows(154).Insert()

# Populate the newly inserted row 154 with the new weekly record
$ws.Cells.Item(154, 1).Value = 5
$ws.Cells.Item(154, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(154, 3).Value = 'Maule'
$ws.Cells.Item(154, 4).Value = 45265
$ws.Cells.Item(154, 4).NumberFormat = $ws.Cells.Item(155, 4).NumberFormat
$ws.Cells.Item(154, 5).Value = 7
$ws.Cells.Item(154, 6).Value = 100112022
$ws.Cells.Item(154, 7).Value = 'Arveja Verde'
$ws.Cells.Item(154, 8).Value = 'Sin especificar'
$ws.Cells.Item(154, 9).Value = 'Primera'
$ws.Cells.Item(154, 10).Value = 500
$ws.Cells.Item(154, 11).Value = 22000
$ws.Cells.Item(154, 12).Value = 22000
$ws.Cells.Item(154, 13).Value = 22000
$ws.Cells.Item(154, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(154, 15).Value = 'Región del Maule'
$ws.Cells.Item(154, 16).Value = 880
$ws.Cells.Item(154, 17).Value = 25
$ws.Cells.Item(154, 18).Value = 'Hortaliza'
